# Weekly update: insert two new price records (rows 89-90) for
# "Cebollín" at Terminal Hortofrutícola Agro Chillán, shifting the
# existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 89:90 - everything from the old row 89
# downward (through the old row 149) shifts down to rows 91-151.
$ws.Rows("89:90").Insert()

# Populate the newly inserted row 89.
$ws.Range("A89").Value = 7
$ws.Range("B89").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C89").Value = 'Ñuble'
$ws.Range("D89").Value = 45086
$ws.Range("E89").Value = 16
$ws.Range("F89").Value = 100112037
$ws.Range("G89").Value = 'Cebollín'
$ws.Range("H89").Value = 'Sin especificar'
$ws.Range("I89").Value = 'Primera'
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = 6000
$ws.Range("N89").Value = '$/paquete 36 unidades'
$ws.Range("O89").Value = 'Provincia de Diguillín'
$ws.Range("P89").Value = 167
$ws.Range("Q89").Value = 36
$ws.Range("R89").Value = 'Hortaliza'

# Populate the newly inserted row 90.
$ws.Range("A90").Value = 7
$ws.Range("B90").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C90").Value = 'Ñuble'
$ws.Range("D90").Value = 45086
$ws.Range("E90").Value = 16
$ws.Range("F90").Value = 100112037
$ws.Range("G90").Value = 'Cebollín'
$ws.Range("H90").Value = 'Sin especificar'
$ws.Range("I90").Value = 'Segunda'
$ws.Range("J90").Value = 80
$ws.Range("K90").Value = 5000
$ws.Range("L90").Value = 5000
$ws.Range("M90").Value = 5000
$ws.Range("N90").Value = '$/paquete 36 unidades'
$ws.Range("O90").Value = 'Provincia de Diguillín'
$ws.Range("P90").Value = 139
$ws.Range("Q90").Value = 36
$ws.Range("R90").Value = 'Hortaliza'
